$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (part CPF0603B75RE1 / Mouser # 279-CPF0603B75RE1) description text
# was corrected from the older "Thin Film Resistors - SMD .1W 75ohm 1% 0603
# 50ppm Auto" wording to the manufacturer's CPF-series spec string.
$ws.Range("C11").Value = "Thin Film Resistors - SMD CPF 0603 75R 0.1% 25PPM "

# Reflect the view state left in the workbook: scrolled one column to the
# right (column B becomes the left-most visible column) with C11 selected.
try {
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}

$ws.Range("C11").Select()
